# Refresh the cryptos price/volume table with the latest scrape (GitHub Actions job).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.623.96"
$ws.Range("E2").Value = "  +0.64%  "
$ws.Range("D3").Value = "2.597.10"
$ws.Range("E3").Value = "  -0.81%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "'590.18"
$ws.Range("E5").Value = "  +0.66%  "
$ws.Range("D6").Value = "'145.94"
$ws.Range("E6").Value = "  -1.88%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  -2.13%  "
$ws.Range("E9").Value = "  -2.75%  "
$ws.Range("E10").Value = "  -0.54%  "
$ws.Range("E12").Value = "  -1.44%  "
$ws.Range("D13").Value = "'27.23"
$ws.Range("E13").Value = "  -1.43%  "
$ws.Range("D14").Value = "3.058.85"
$ws.Range("E14").Value = "  -0.87%  "
$ws.Range("D15").Value = "63.350.72"
$ws.Range("E15").Value = "  +0.43%  "
$ws.Range("E16").Value = "  -1.83%  "
$ws.Range("D17").Value = "2.588.55"
$ws.Range("E17").Value = "  -3.73%  "
$ws.Range("D18").Value = "'11.17"
$ws.Range("E18").Value = "  -2.61%  "
$ws.Range("D19").Value = "'342.36"
$ws.Range("E19").Value = "  -1.18%  "
$ws.Range("D20").Value = "'4.34"
$ws.Range("E20").Value = "  -2.51%  "
$ws.Range("D21").Value = "'6.70"
$ws.Range("E21").Value = "  -2.27%  "
$ws.Range("E22").Value = "  -0.08%  "
$ws.Range("D23").Value = "'68.43"
$ws.Range("E23").Value = "  +2.20%  "
$ws.Range("D24").Value = "'1.57"
$ws.Range("E24").Value = "  +5.75%  "
$ws.Range("D25").Value = "'1.62"
$ws.Range("E25").Value = "  +0.46%  "
$ws.Range("D26").Value = "'0.166"
$ws.Range("E26").Value = "  -3.09%  "
$ws.Range("D27").Value = "'0.997"
$ws.Range("E27").Value = "  -0.70%  "
$ws.Range("D28").Value = "'7.86"
$ws.Range("E28").Value = "  -3.62%  "
$ws.Range("D29").Value = "'8.30"
$ws.Range("E29").Value = "  -2.48%  "
$ws.Range("D30").Value = "'1.97"
$ws.Range("E30").Value = "  -0.44%  "
$ws.Range("D31").Value = "'478.21"
$ws.Range("E31").Value = "  +2.08%  "
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").Value = "'1.72"
$ws.Range("E32").Value = "  +5.73%  "
$ws.Range("B33").Value = "PEPE"
$ws.Range("C33").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D33").Value = "0.0₃0813"
$ws.Range("E33").Value = "  -2.93%  "
$ws.Range("D34").Value = "'176.42"
$ws.Range("E34").Value = "  +0.09%  "
$ws.Range("E35").Value = "  +0.15%  "
$ws.Range("D36").Value = "'0.397"
$ws.Range("E36").Value = "  -2.80%  "
$ws.Range("D37").Value = "'18.92"
$ws.Range("E37").Value = "  -2.43%  "
$ws.Range("D38").Value = "'4.54"
$ws.Range("E38").Value = "  -1.43%  "
$ws.Range("D39").Value = "'0.998"
$ws.Range("E39").Value = "  -0.13%  "
$ws.Range("E40").Value = "  -0.15%  "
$ws.Range("D41").Value = "'162.87"
$ws.Range("E41").Value = "  +4.14%  "
$ws.Range("D42").Value = "'40.17"
$ws.Range("E42").Value = "  +1.32%  "
$ws.Range("D43").Value = "'3.74"
$ws.Range("E43").Value = "  -2.51%  "
$ws.Range("D44").Value = "'21.76"
$ws.Range("E44").Value = "  +3.72%  "
$ws.Range("D45").Value = "'0.627"
$ws.Range("E45").Value = "  -3.14%  "
$ws.Range("D46").Value = "'0.0539"
$ws.Range("E46").Value = "  -2.58%  "
$ws.Range("D47").Value = "'0.0964"
$ws.Range("E47").Value = "  -1.30%  "
$ws.Range("D48").Value = "'0.0238"
$ws.Range("E48").Value = "  -1.23%  "
$ws.Range("D49").Value = "'18.46"
$ws.Range("E49").Value = "  -3.11%  "
$ws.Range("B50").Value = "dogwifhat"
$ws.Range("C50").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D50").Value = "'1.74"
$ws.Range("E50").Value = "  -0.78%  "
$ws.Range("B51").Value = "WhiteBITCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D51").Value = "'11.37"
$ws.Range("E51").Value = "  -0.38%  "
